$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties, reusing the existing header style (s="1")
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record values (Wins=76, Losses=86, Ties=0) for all data rows (2-60)
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD
    $ws.Cells.Item($r, 31).Value = 86   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
